$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Update GAME NAME / GAME NUMBER / TOP PRIZES REMAINING for rows that shifted ---
$ws.Cells.Item(2, 3).Value = "Super 7s Jackpot"
$ws.Cells.Item(2, 4).Value = 2117
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(3, 3).Value = "Lucky 8s Doubler"
$ws.Cells.Item(3, 4).Value = 2091
$ws.Cells.Item(3, 5).Value = 6
$ws.Cells.Item(5, 3).Value = "Find The 9s"
$ws.Cells.Item(5, 4).Value = 2130
$ws.Cells.Item(5, 5).Value = 239
$ws.Cells.Item(9, 3).Value = "Lucky Dog"
$ws.Cells.Item(9, 4).Value = 2071
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(14, 3).Value = "Money Multipler"
$ws.Cells.Item(14, 4).Value = 2108
$ws.Cells.Item(14, 5).Value = 4
$ws.Cells.Item(16, 3).Value = "Weekly Grand"
$ws.Cells.Item(16, 4).Value = 2078
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(17, 3).Value = "$30,000 Jackpot"
$ws.Cells.Item(17, 4).Value = 2133
$ws.Cells.Item(17, 5).Value = 8
$ws.Cells.Item(23, 3).Value = "777"
$ws.Cells.Item(23, 4).Value = 2035
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(25, 3).Value = "Cinco Connect"
$ws.Cells.Item(25, 4).Value = 2031
$ws.Cells.Item(25, 5).Value = 4
$ws.Cells.Item(26, 3).Value = "Texas Loteria"
$ws.Cells.Item(26, 4).Value = 1828
$ws.Cells.Item(26, 5).Value = 22
$ws.Cells.Item(27, 3).Value = "$50,000 Bonus Cashword"
$ws.Cells.Item(27, 4).Value = 2052
$ws.Cells.Item(27, 5).Value = 17
$ws.Cells.Item(30, 3).Value = "Diamond Mine 9X"
$ws.Cells.Item(30, 4).Value = 2090
$ws.Cells.Item(30, 5).Value = 1
$ws.Cells.Item(31, 3).Value = "Super Loteria"
$ws.Cells.Item(31, 4).Value = 1877
$ws.Cells.Item(31, 5).Value = 2
$ws.Cells.Item(32, 3).Value = "Frogger"
$ws.Cells.Item(32, 4).Value = 2049
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(33, 3).Value = "Multiplier Mania"
$ws.Cells.Item(33, 4).Value = 2106
$ws.Cells.Item(33, 5).Value = 3
$ws.Cells.Item(34, 3).Value = "Bonus Break the Bank"
$ws.Cells.Item(34, 4).Value = 1862
$ws.Cells.Item(34, 5).Value = 2
$ws.Cells.Item(36, 3).Value = "Golden Nugget"
$ws.Cells.Item(36, 4).Value = 2116
$ws.Cells.Item(36, 5).Value = 4
$ws.Cells.Item(37, 3).Value = "Cowboys"
$ws.Cells.Item(37, 4).Value = 2069
$ws.Cells.Item(37, 5).Value = 7
$ws.Cells.Item(38, 3).Value = "Money Madness"
$ws.Cells.Item(38, 4).Value = 2102
$ws.Cells.Item(38, 5).Value = 3
$ws.Cells.Item(42, 3).Value = "Día De Los Muertos"
$ws.Cells.Item(42, 4).Value = 2092
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(43, 3).Value = "Crazy Cash Blowout"
$ws.Cells.Item(43, 4).Value = 2125
$ws.Cells.Item(43, 5).Value = 4
$ws.Cells.Item(44, 3).Value = "Super Loteria"
$ws.Cells.Item(44, 4).Value = 2074
$ws.Cells.Item(44, 5).Value = 18
$ws.Cells.Item(46, 3).Value = "Lucky 7s Hunt"
$ws.Cells.Item(46, 4).Value = 2055
$ws.Cells.Item(46, 5).Value = 2
$ws.Cells.Item(47, 3).Value = "Super Crossword"
$ws.Cells.Item(47, 4).Value = 2082
$ws.Cells.Item(47, 5).Value = 9
$ws.Cells.Item(53, 3).Value = "Red Hot 7s"
$ws.Cells.Item(53, 4).Value = 2095
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 3).Value = "50X Fast Cash"
$ws.Cells.Item(54, 4).Value = 2063
$ws.Cells.Item(54, 5).Value = 0
$ws.Cells.Item(56, 3).Value = "$250,000 50X Cashword"
$ws.Cells.Item(56, 4).Value = 2083
$ws.Cells.Item(56, 5).Value = 2
$ws.Cells.Item(58, 3).Value = "Casino Action Super Ticket™"
$ws.Cells.Item(58, 4).Value = 2076
$ws.Cells.Item(58, 5).Value = 2
$ws.Cells.Item(62, 3).Value = "Cash Blast"
$ws.Cells.Item(62, 4).Value = 2103
$ws.Cells.Item(62, 5).Value = 896
$ws.Cells.Item(72, 3).Value = "$100 or $200"
$ws.Cells.Item(72, 4).Value = 2093
$ws.Cells.Item(72, 5).Value = 128186

# --- Step 2: Update LAST SCRAPE DATE column, keeping values as text (not auto-converted dates) ---
# Build helper cells with the exact date strings as text using a formula-then-paste-values trick,
# which avoids Excel auto-converting the text into a date serial number.
$ws.Cells.Item(1, 8).Formula = "=""2019-02-10"""
$ws.Cells.Item(2, 8).Formula = "=""2019-02-11"""
$ws.Cells.Item(3, 8).Formula = "=""2019-02-27"""
$ws.Cells.Item(4, 8).Formula = "=""2019-03-12"""

# Copy each helper text value into the destination date cells via PasteSpecial (values only)
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(54, 6).PasteSpecial(-4163)
$ws.Cells.Item(2, 8).Copy()
$ws.Cells.Item(53, 6).PasteSpecial(-4163)
$ws.Cells.Item(3, 8).Copy()
$ws.Cells.Item(34, 6).PasteSpecial(-4163)
$ws.Cells.Item(4, 8).Copy()
$ws.Cells.Item(2, 6).PasteSpecial(-4163)
$ws.Cells.Item(3, 6).PasteSpecial(-4163)
$ws.Cells.Item(4, 6).PasteSpecial(-4163)
$ws.Cells.Item(5, 6).PasteSpecial(-4163)
$ws.Cells.Item(6, 6).PasteSpecial(-4163)
$ws.Cells.Item(7, 6).PasteSpecial(-4163)
$ws.Cells.Item(8, 6).PasteSpecial(-4163)
$ws.Cells.Item(9, 6).PasteSpecial(-4163)
$ws.Cells.Item(12, 6).PasteSpecial(-4163)
$ws.Cells.Item(13, 6).PasteSpecial(-4163)
$ws.Cells.Item(14, 6).PasteSpecial(-4163)
$ws.Cells.Item(15, 6).PasteSpecial(-4163)
$ws.Cells.Item(16, 6).PasteSpecial(-4163)
$ws.Cells.Item(17, 6).PasteSpecial(-4163)
$ws.Cells.Item(18, 6).PasteSpecial(-4163)
$ws.Cells.Item(19, 6).PasteSpecial(-4163)
$ws.Cells.Item(20, 6).PasteSpecial(-4163)
$ws.Cells.Item(21, 6).PasteSpecial(-4163)
$ws.Cells.Item(23, 6).PasteSpecial(-4163)
$ws.Cells.Item(24, 6).PasteSpecial(-4163)
$ws.Cells.Item(25, 6).PasteSpecial(-4163)
$ws.Cells.Item(26, 6).PasteSpecial(-4163)
$ws.Cells.Item(27, 6).PasteSpecial(-4163)
$ws.Cells.Item(28, 6).PasteSpecial(-4163)
$ws.Cells.Item(29, 6).PasteSpecial(-4163)
$ws.Cells.Item(30, 6).PasteSpecial(-4163)
$ws.Cells.Item(31, 6).PasteSpecial(-4163)
$ws.Cells.Item(32, 6).PasteSpecial(-4163)
$ws.Cells.Item(33, 6).PasteSpecial(-4163)
$ws.Cells.Item(35, 6).PasteSpecial(-4163)
$ws.Cells.Item(36, 6).PasteSpecial(-4163)
$ws.Cells.Item(37, 6).PasteSpecial(-4163)
$ws.Cells.Item(38, 6).PasteSpecial(-4163)
$ws.Cells.Item(39, 6).PasteSpecial(-4163)
$ws.Cells.Item(40, 6).PasteSpecial(-4163)
$ws.Cells.Item(41, 6).PasteSpecial(-4163)
$ws.Cells.Item(42, 6).PasteSpecial(-4163)
$ws.Cells.Item(43, 6).PasteSpecial(-4163)
$ws.Cells.Item(44, 6).PasteSpecial(-4163)
$ws.Cells.Item(45, 6).PasteSpecial(-4163)
$ws.Cells.Item(46, 6).PasteSpecial(-4163)
$ws.Cells.Item(47, 6).PasteSpecial(-4163)
$ws.Cells.Item(48, 6).PasteSpecial(-4163)
$ws.Cells.Item(49, 6).PasteSpecial(-4163)
$ws.Cells.Item(50, 6).PasteSpecial(-4163)
$ws.Cells.Item(51, 6).PasteSpecial(-4163)
$ws.Cells.Item(55, 6).PasteSpecial(-4163)
$ws.Cells.Item(56, 6).PasteSpecial(-4163)
$ws.Cells.Item(57, 6).PasteSpecial(-4163)
$ws.Cells.Item(58, 6).PasteSpecial(-4163)
$ws.Cells.Item(59, 6).PasteSpecial(-4163)
$ws.Cells.Item(60, 6).PasteSpecial(-4163)
$ws.Cells.Item(61, 6).PasteSpecial(-4163)
$ws.Cells.Item(62, 6).PasteSpecial(-4163)
$ws.Cells.Item(63, 6).PasteSpecial(-4163)
$ws.Cells.Item(64, 6).PasteSpecial(-4163)
$ws.Cells.Item(65, 6).PasteSpecial(-4163)
$ws.Cells.Item(66, 6).PasteSpecial(-4163)
$ws.Cells.Item(67, 6).PasteSpecial(-4163)
$ws.Cells.Item(68, 6).PasteSpecial(-4163)
$ws.Cells.Item(69, 6).PasteSpecial(-4163)
$ws.Cells.Item(70, 6).PasteSpecial(-4163)
$ws.Cells.Item(71, 6).PasteSpecial(-4163)
$ws.Cells.Item(72, 6).PasteSpecial(-4163)
$ws.Cells.Item(73, 6).PasteSpecial(-4163)
$ws.Cells.Item(74, 6).PasteSpecial(-4163)
$ws.Cells.Item(75, 6).PasteSpecial(-4163)
$ws.Cells.Item(76, 6).PasteSpecial(-4163)
$ws.Cells.Item(77, 6).PasteSpecial(-4163)
$ws.Cells.Item(78, 6).PasteSpecial(-4163)
$ws.Cells.Item(79, 6).PasteSpecial(-4163)
$ws.Cells.Item(80, 6).PasteSpecial(-4163)
$ws.Cells.Item(81, 6).PasteSpecial(-4163)

# --- Step 3: Clean up helper cells ---
$ws.Range($ws.Cells.Item(1,8), $ws.Cells.Item(4,8)).ClearContents()
$excel.CutCopyMode = 0
